# DOMA-1856 add account number column in export meter readings
#
# Inserts a new "Лицевой счет" (account number) column between the existing
# "Квартира" (C) and "Услуга" (D) columns, shifting columns D:L to E:M, and
# fills the header + the two template placeholder rows for the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing D:L columns one position to the right, opening up a
# fresh column D for the new "account number" field.
$ws.Range("D1:D3").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "Лицевой счет"

# Template rows (merge-print placeholders consumed by the export engine)
$ws.Range("D2").Value = "{d.meter[i].accountNumber}"
$ws.Range("D3").Value = "{d.meter[i + 1].accountNumber}"
